$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 119; this shifts the existing rows 119-126
# down to 120-127, matching the rest of the diff (those rows keep their
# original values, just moved down by one row).
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new record's data.
$ws.Range("A119").Value = 6
$ws.Range("B119").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C119").Value = "Metropolitana"
$ws.Range("D119").Value = 44516
$ws.Range("E119").Value = 13
$ws.Range("F119").Value = 100112001
$ws.Range("G119").Value = "Berenjena"
$ws.Range("H119").Value = "Sin especificar"
$ws.Range("I119").Value = "Primera"
$ws.Range("J119").Value = 200
$ws.Range("K119").Value = 13000
$ws.Range("L119").Value = 15000
$ws.Range("M119").Value = 13800
$ws.Range("N119").Value = "$/caja 60 unidades"
$ws.Range("O119").Value = "Provincia de Huasco"
$ws.Range("P119").Value = 230
$ws.Range("Q119").Value = 60
$ws.Range("R119").Value = "Hortaliza"
